{"js": "// 1) Merge the \"PL1.1\" + \" create product line with correct input ext rep\"\n//    runs in the PL1.1 test-case paragraph into a single run.\n// 2) Add three new test-case paragraphs (ITM1.4, ITM1.5, ITM1.6) right\n//    after the \"ITM1.3 create item with bad data\" paragraph, matching its\n//    formatting.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet pl11Paragraph = null;\nlet itm13Paragraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (pl11Paragraph === null && text.indexOf(\"PL1.1\") === 0) {\n    pl11Paragraph = paragraphs.items[i];\n  }\n  if (itm13Paragraph === null && text.indexOf(\"ITM1.3 create item with bad data\") === 0) {\n    itm13Paragraph = paragraphs.items[i];\n  }\n}\n\nif (pl11Paragraph) {\n  // Rewriting the paragraph text collapses the two runs (\"PL1.1\" and\n  // \" create product line with correct input ext rep\") into one run.\n  pl11Paragraph.insertText(\"PL1.1 create product line with correct input ext rep\", \"Replace\");\n}\n\nif (itm13Paragraph) {\n  // insertParagraph inherits the anchor paragraph's formatting (pPr/rPr),\n  // so each new paragraph keeps ITM1.3's underline/Arial styling.\n  let anchor = itm13Paragraph;\n  const newTexts = [\n    \"ITM1.4 create new Department\",\n    \"ITM1.5 create new Category\",\n    \"ITM1.6 create new Subcategory\"\n  ];\n  for (const t of newTexts) {\n    anchor = anchor.insertParagraph(t, \"After\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# 1) Merge the \"PL1.1\" + \" create product line with correct input ext rep\"\n#    runs in the PL1.1 test-case paragraph into a single run.\n# 2) Add three new test-case paragraphs (ITM1.4, ITM1.5, ITM1.6) right\n#    after the \"ITM1.3 create item with bad data\" paragraph, matching its\n#    formatting.\n\n$d = $word.ActiveDocument\n\n# --- Part 1: merge the PL1.1 paragraph's runs ---\n$pl11 = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"PL1.1\")) {\n        $pl11 = $p\n        break\n    }\n}\n\nif ($pl11 -ne $null) {\n    $pr = $pl11.Range\n    # Exclude the trailing paragraph mark from the range used for Find/Replace.\n    $textRange = $d.Range($pr.Start, $pr.End - 1)\n    $fullText = $textRange.Text\n    # Replacing the range's text with itself via Find/Replace rewrites the\n    # paragraph's runs as a single run (collapsing the original two runs)\n    # while keeping the existing run formatting.\n    $textRange.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, $fullText, 2) | Out-Null\n}\n\n# --- Part 2: insert the three new test-case paragraphs after ITM1.3 ---\n$itm13Index = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"ITM1.3 create item with bad data\")) {\n        $itm13Index = $i\n        break\n    }\n}\n\nif ($itm13Index -gt 0) {\n    $newTexts = @(\"ITM1.4 create new Department\", \"ITM1.5 create new Category\", \"ITM1.6 create new Subcategory\")\n    $curIndex = $itm13Index\n    foreach ($t in $newTexts) {\n        $anchor = $d.Paragraphs.Item($curIndex)\n        # InsertParagraphAfter creates a new paragraph that inherits the\n        # anchor paragraph's formatting (pPr/rPr), matching ITM1.3's style.\n        $anchor.Range.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($curIndex + 1)\n        $newPara.Range.InsertBefore($t)\n        $curIndex = $curIndex + 1\n    }\n}\n"}
